$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [char]0x2083

$ws.Range("D2").Value = "59.034.05"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "2.508.95"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'536.57"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").Value = "'137.31"
$ws.Range("E6").Value = "  -2.22%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").Value = "'0.562"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").Value = "2.535.84"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("E10").Value = "  +2.24%  "

$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").Value = "'5.31"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("D13").Value = "'0.349"
$ws.Range("E13").Value = "  -0.96%  "

$ws.Range("D14").Value = "2.951.59"
$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").Value = "'23.22"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").Value = "59.168.07"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").Value = "2.520.11"
$ws.Range("E18").Value = "  -1.53%  "

$ws.Range("D19").Value = "'11.08"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").Value = "'4.27"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").Value = "'326.36"
$ws.Range("E21").Value = "  +2.01%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("D24").Value = "'63.48"
$ws.Range("E24").Value = "  +1.80%  "

$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").Value = "'7.60"
$ws.Range("E28").Value = "  -2.56%  "

$ws.Range("E29").Value = "  +1.05%  "

$ws.Range("D30").Value = "0.0" + $sub3.ToString() + "0779"
$ws.Range("E30").Value = "  +1.66%  "

$ws.Range("D31").Value = "'1.78"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").Value = "'166.67"
$ws.Range("E32").Value = "  +1.46%  "

$ws.Range("D33").Value = "'1.15"
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").Value = "'1.41"
$ws.Range("E35").Value = "  -1.71%  "

$ws.Range("D36").Value = "'18.51"
$ws.Range("E36").Value = "  +0.77%  "

$ws.Range("D37").Value = "'4.14"
$ws.Range("E37").Value = "  -2.71%  "

$ws.Range("D38").Value = "'1.56"
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D39").Value = "'36.79"
$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("D40").Value = "'0.831"
$ws.Range("E40").Value = "  +3.45%  "

$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").Value = "'5.25"
$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").Value = "'278.99"
$ws.Range("E43").Value = "  -2.98%  "

$ws.Range("D44").Value = "'0.993"
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("E45").Value = "  +1.61%  "

$ws.Range("D46").Value = "'10.83"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("D47").Value = "'125.69"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").Value = "'0.0222"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").Value = "'17.61"
$ws.Range("E51").Value = "  +0.99%  "
